$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', None), ('selector', None),
                ('model',
                 BaggingClassifier(estimator=SVC(C=5, class_weight='balanced',
                                                 kernel='poly',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])"
$ws.Range("B2").Value = 0.6571428571428571
$ws.Range("C2").Value = "{'selector': None, 'scaler': None, 'model__n_estimators': 50, 'model__estimator__kernel': 'poly', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 5}"
$ws.Range("D2").Value = 0.5333333333333333
$ws.Range("E2").Value = "[1 0 0 1 0 0 1 1 0 1 0 0]"
$ws.Range("F2").Value = "[1 0 1 0 1 1 1 1 1 1 1 1]"
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.670452380952381
$ws.Range("I2").Value = 0.02891317097757153
$ws.Range("J2").Value = 0.5587619047619047
$ws.Range("K2").Value = 0.05555168013123375

# Row 3
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', None), ('selector', None),
                ('model',
                 BaggingClassifier(estimator=SVC(C=1, class_weight='balanced',
                                                 random_state=42),
                                   n_estimators=50, random_state=42))])"
$ws.Range("B3").Value = 0.6476190476190476
$ws.Range("C3").Value = "{'selector': None, 'scaler': None, 'model__n_estimators': 50, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': 'balanced', 'model__estimator__C': 1}"
$ws.Range("D3").Value = 0.5
$ws.Range("E3").Value = "[1 0 1 0 0 0 0 1 1 0 1 1]"
$ws.Range("F3").Value = "[1 1 1 1 1 0 0 0 0 0 1 0]"
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.6895238095238093
$ws.Range("I3").Value = 0.03309399110085322
$ws.Range("J3").Value = 0.5476190476190476
$ws.Range("K3").Value = 0.05641281351300525

# Row 4
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', None),
                ('selector',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                ('model',
                 BaggingClassifier(estimator=SVC(C=1, random_state=42),
                                   random_state=42))])"
$ws.Range("B4").Value = 0.6380952380952382
$ws.Range("C4").Value = "{'selector': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), 'scaler': None, 'model__n_estimators': 10, 'model__estimator__kernel': 'rbf', 'model__estimator__class_weight': None, 'model__estimator__C': 1}"
$ws.Range("D4").Value = 0.7777777777777777
$ws.Range("E4").Value = "[1 0 1 1 1 1 0 1 0 1 0 1]"
$ws.Range("F4").Value = "[1 1 1 1 1 1 0 0 1 1 1 1]"
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6738095238095237
$ws.Range("I4").Value = 0.03559195631867585
$ws.Range("J4").Value = 0.5260952380952381
$ws.Range("K4").Value = 0.06835069464204943
